$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Time Format" row beneath the existing data (row 5)
$ws.Range("A5").Value = "Time Format"
$ws.Range("B5").Value = 0.12934027777777779
$ws.Range("B5").NumberFormat = "h:mm:ss"

# Move / update the active selection to match the saved view state
$ws.Range("F17").Select()
